# The diff shows rows 5-9 being cyclically re-shuffled: the *content* that
# used to live in row 8 now lives in row 5, row 9's content now lives in
# row 6, row 5's content now lives in row 7, row 6's content now lives in
# row 8, and row 7's content now lives in row 9. (Matching Id/TaxonId/
# Artnamn values identify this as a pure row permutation, not per-cell
# edits.)
#
# sourceForTarget[targetRow] = rowThatSuppliesTheNewContent
$sourceForTarget = @{ 5 = 8; 6 = 9; 7 = 5; 8 = 6; 9 = 7 }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 51  # AY

# Columns I (9) and AR (44) hold values that are sometimes purely numeric
# looking text (e.g. "10", "16317") which Excel would otherwise
# auto-convert to a real number when assigned via Value2. Force those
# through as text so the stored type stays text, same as the source cells.
$textCols = @(9, 44)

# Columns Y/Z/AA/AB (Startdatum/Starttid/Slutdatum/Sluttid) hold date-like
# text ("2013-09-17"/"00:00") that is identical across every one of rows
# 5-9 in the source data, AND which Excel would otherwise auto-convert to
# a real date/time serial number when re-assigned via Value2. Since the
# rotation never actually changes these columns' values, just skip
# rewriting them entirely and leave the original cells untouched.
$skipCols = @(25, 26, 27, 28)

$rows = @(5, 6, 7, 8, 9)

# 1) Snapshot every cell (A:AY) of the affected rows before any writes,
#    so that reads never see already-overwritten data.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowData[$col] = $ws.Cells.Item($r, $col).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write each target row's cells from the snapshot of its source row.
foreach ($r in $rows) {
    $srcRow = $sourceForTarget[$r]
    $rowData = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($skipCols -contains $col) {
            continue
        }
        $val = $rowData[$col]
        $cell = $ws.Cells.Item($r, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value2 = $val
            $cell.ClearFormats()
        } else {
            $cell.Value2 = $val
        }
    }
}
